$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 12 values/styles (G12, H12, I12, J12) ---
# G12 and H12 already carry style index 2; just fill in the values.
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 5

# I12 needs the style currently used by J15/J21 (border-only style, index 8)
# and J12 needs the style currently used by J9 (fill/border style, index 9).
# Copy formats from those reference cells, then set the values.
$ws.Range("J15").Copy()
$ws.Range("I12").PasteSpecial(-4122)

$ws.Range("J9").Copy()
$ws.Range("J12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 5

# --- Update the frozen-pane view / active selection ---
# Re-establish the freeze at the same boundary (2 columns x 3 rows) and move
# the active selection from J14 to J12, matching the updated sheet view.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("C4").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("J12").Select()
